# Adds the newly-scraped pharmacy leads (rows 223-227) to the sheet,
# matching the existing "İşletme Adı / Telefon / Adres / Bölge / Tarih" layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(223,1).Value = "ECZANE Vatan"
$ws.Cells.Item(223,2).Value = "
+90 216 565 96 71"
$ws.Cells.Item(223,3).Value = "
Dumlupınar Mh.şahika sok.no:1/C FİKİRTEPE, 34720 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(223,4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(223,5).Value = "2026-02-12 01:27"
$ws.Rows.Item(223).AutoFit()

$ws.Cells.Item(224,1).Value = "Pelikan Eczanesi"
$ws.Cells.Item(224,2).Value = "
+90 530 765 82 94"
$ws.Cells.Item(224,3).Value = "
Göztepe, Tütüncü Mehmet Efendi Cd. 107/B, 34730 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(224,4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(224,5).Value = "2026-02-12 01:27"
$ws.Rows.Item(224).AutoFit()

$ws.Cells.Item(225,1).Value = "Ecem Ramiz Eczanesi"
$ws.Cells.Item(225,2).Value = "
+90 216 337 09 19"
$ws.Cells.Item(225,3).Value = "
Göztepe, Nadirağa Sk. No:5-7E D:F, 34730 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(225,4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(225,5).Value = "2026-02-12 01:27"
$ws.Rows.Item(225).AutoFit()

$ws.Cells.Item(226,1).Value = "Caddebostan Plus Eczanesi"
$ws.Cells.Item(226,2).Value = "
+90 545 350 10 00"
$ws.Cells.Item(226,3).Value = "
Caddebostan, Ömer Paşa Sk. no:1, 34728 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(226,4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(226,5).Value = "2026-02-12 01:27"
$ws.Rows.Item(226).AutoFit()

$ws.Cells.Item(227,1).Value = "Kekik Eczanesi"
$ws.Cells.Item(227,2).Value = "
+90 216 356 78 91"
$ws.Cells.Item(227,3).Value = "
Atatürk Caddesi Makbule Apt. No.37/A Sahrayıcedit, 34734 İstanbul, Türkiye"
$ws.Cells.Item(227,4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(227,5).Value = "2026-02-12 01:27"
$ws.Rows.Item(227).AutoFit()

